$wb = $excel.ActiveWorkbook

# ---- Sheet "Jorge Luis" (sheet1): row 12 ----
$ws1 = $wb.Worksheets.Item("Jorge Luis")

$ws1.Cells.Item(12, 1).Value = 41704
$ws1.Cells.Item(11, 1).Copy()
$ws1.Cells.Item(12, 1).PasteSpecial(-4122)
$ws1.Cells.Item(12, 2).Value = "crear interfaz de el pedido y carrito"
$ws1.Cells.Item(12, 3).Value = 3
$ws1.Cells.Item(12, 6).Value = 70

# ---- Sheet "Fabio" (sheet2): rows 11 and 12 ----
$ws2 = $wb.Worksheets.Item("Fabio")

# Row 11
$ws2.Cells.Item(11, 1).Value = 41704
$ws2.Cells.Item(10, 1).Copy()
$ws2.Cells.Item(11, 1).PasteSpecial(-4122)
$ws2.Cells.Item(11, 2).Value = "Imágenes a utilizar en portafolio"
$ws2.Cells.Item(11, 4).Value = 16
$ws2.Cells.Item(11, 6).Value = 100

# Row 12 (introduces new shared string "terminar portafolio" -> index 38)
$ws2.Cells.Item(12, 1).Value = 41704
$ws2.Cells.Item(10, 1).Copy()
$ws2.Cells.Item(12, 1).PasteSpecial(-4122)
$ws2.Cells.Item(12, 2).Value = "terminar portafolio"
$ws2.Cells.Item(12, 3).Value = 5

# ---- Sheet "Jorge Luis" (sheet1): row 13 (introduces new shared string "agregar el pedido a una base de datos" -> index 39) ----
$ws1.Cells.Item(13, 1).Value = 41704
$ws1.Cells.Item(11, 1).Copy()
$ws1.Cells.Item(13, 1).PasteSpecial(-4122)
$ws1.Cells.Item(13, 2).Value = "agregar el pedido a una base de datos"
$ws1.Cells.Item(13, 3).Value = 8
$ws1.Cells.Item(13, 6).Value = 0

$excel.CutCopyMode = 0

# Move the active selection on "Fabio" down to the newly added last row,
# then restore "Jorge Luis" as the active/selected sheet tab.
$ws2.Range("A12").Select() | Out-Null
$ws1.Activate() | Out-Null

